# Update cryptocurrency price/volume data (refresh from GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell holding the default/unstyled "Normal" style, used to restore
# formatting on cells after a temporary text-number-format is applied so that
# numeric-looking price strings (e.g. "43.00") are stored as text, not numbers.
$normalStyle = $ws.Range("Z100").Style

# Row 2
$ws.Range("D2").Value = "65.202.11"
$ws.Range("E2").Value = "  -0.66%  "

# Row 3
$ws.Range("D3").Value = "3.537.41"
$ws.Range("E3").Value = "  +2.77%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.38"
$ws.Range("D5").Style = $normalStyle
$ws.Range("E5").Value = "  +1.17%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.99"
$ws.Range("D6").Style = $normalStyle
$ws.Range("E6").Value = "  +1.24%  "

# Row 7
$ws.Range("D7").Value = "3.537.62"
$ws.Range("E7").Value = "  +2.80%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.489"
$ws.Range("D9").Style = $normalStyle
$ws.Range("E9").Value = "  -2.21%  "

# Row 10
$ws.Range("E10").Value = "  +2.03%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.91"
$ws.Range("D11").Style = $normalStyle
$ws.Range("E11").Value = "  -5.91%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.391"
$ws.Range("D12").Style = $normalStyle
$ws.Range("E12").Value = "  +3.03%  "

# Row 13
$ws.Range("D13").Value = "4.137.55"
$ws.Range("E13").Value = "  +2.88%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000186"
$ws.Range("D14").Style = $normalStyle
$ws.Range("E14").Value = "  +1.84%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.12"
$ws.Range("D15").Style = $normalStyle
$ws.Range("E15").Value = "  +2.56%  "

# Row 16
$ws.Range("D16").Value = "3.538.25"
$ws.Range("E16").Value = "  +3.01%  "

# Row 17
$ws.Range("E17").Value = "  +1.55%  "

# Row 18
$ws.Range("D18").Value = "65.293.61"
$ws.Range("E18").Value = "  -0.43%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.31"
$ws.Range("D19").Style = $normalStyle
$ws.Range("E19").Value = "  +4.65%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.94"
$ws.Range("D20").Style = $normalStyle
$ws.Range("E20").Value = "  +0.94%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.34"
$ws.Range("D21").Style = $normalStyle
$ws.Range("E21").Value = "  +4.32%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "393.31"
$ws.Range("D22").Style = $normalStyle
$ws.Range("E22").Value = "  -0.12%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.574"
$ws.Range("D23").Style = $normalStyle
$ws.Range("E23").Value = "  +3.55%  "

# Row 24
$ws.Range("D24").Value = "3.677.05"
$ws.Range("E24").Value = "  +2.62%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.91"
$ws.Range("D25").Style = $normalStyle
$ws.Range("E25").Value = "  +0.72%  "

# Row 26
$ws.Range("E26").Value = "  -0.06%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000116"
$ws.Range("D27").Style = $normalStyle
$ws.Range("E27").Value = "  +7.69%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.75"
$ws.Range("D28").Style = $normalStyle
$ws.Range("E28").Value = "  +7.76%  "

# Row 29
$ws.Range("E29").Value = "  +0.37%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.29"
$ws.Range("D30").Style = $normalStyle
$ws.Range("E30").Value = "  +1.96%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.23"
$ws.Range("D31").Style = $normalStyle
$ws.Range("E31").Value = "  -0.11%  "

# Row 32
$ws.Range("D32").Value = "3.549.21"
$ws.Range("E32").Value = "  +2.92%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.81"
$ws.Range("D34").Style = $normalStyle
$ws.Range("E34").Value = "  +3.50%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.145"
$ws.Range("D35").Style = $normalStyle
$ws.Range("E35").Value = "  -1.02%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.28"
$ws.Range("D36").Style = $normalStyle
$ws.Range("E36").Value = "  +8.19%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.96"
$ws.Range("D37").Style = $normalStyle
$ws.Range("E37").Value = "  +0.22%  "

# Row 38
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "168.76"
$ws.Range("D38").Style = $normalStyle
$ws.Range("E38").Value = "  -2.15%  "

# Row 39
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.55"
$ws.Range("D39").Style = $normalStyle
$ws.Range("E39").Value = "  +4.28%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.95"
$ws.Range("D40").Style = $normalStyle
$ws.Range("E40").Value = "  +2.88%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0810"
$ws.Range("D41").Style = $normalStyle
$ws.Range("E41").Value = "  +5.38%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.826"
$ws.Range("D42").Style = $normalStyle
$ws.Range("E42").Value = "  -0.02%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.54"
$ws.Range("D43").Style = $normalStyle
$ws.Range("E43").Value = "  +15.71%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.00"
$ws.Range("D44").Style = $normalStyle
$ws.Range("E44").Value = "  -1.78%  "

# Row 45
$ws.Range("E45").Value = "  +0.01%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.43"
$ws.Range("D46").Style = $normalStyle
$ws.Range("E46").Value = "  +0.43%  "

# Row 47
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.67"
$ws.Range("D47").Style = $normalStyle
$ws.Range("E47").Value = "  +3.19%  "

# Row 48
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.18"
$ws.Range("D48").Style = $normalStyle
$ws.Range("E48").Value = "  +6.87%  "

# Row 49
$ws.Range("D49").Value = "2.437.88"
$ws.Range("E49").Value = "  +10.13%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.80"
$ws.Range("D50").Style = $normalStyle
$ws.Range("E50").Value = "  +3.14%  "

# Row 51
$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.35"
$ws.Range("D51").Style = $normalStyle
$ws.Range("E51").Value = "  +16.45%  "
